$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 30) with the language code "en" in column A,
# matching the existing "language" column values used throughout the sheet.
$ws.Range("A30").Value = "en"

# Move the active selection to reflect the new bottom of the used range,
# similar to the author's recorded selection after the edit.
$ws.Range("C31").Select()
